# Refresh the crypto price table pulled in by the scheduled GitHub Actions job.
# Rows 2-51 hold one coin each (A=rank, B=Coin, C=Link, D=Price, E=Volume(1h)).
# This run mostly updates the Price / Volume(1h) text for every coin; rows 41 and
# 42 (Quant / TrustWalletToken) additionally swapped ranking order, so their Coin,
# Link, Price and Volume(1h) values all move between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of these columns are stored as plain text in the workbook (e.g. prices such
# as "1.002" or "29.316.18" are text, not numbers). Writing straight to .Value lets
# Excel silently reinterpret anything that looks numeric as a real number, so force
# the cell to text format first and drop back to the default "Normal" style once the
# text is in place, leaving no left-over custom number format on the cell.
function Set-TextValue($Cell, [string]$Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) '29.316.18'  # D2 (Price)
Set-TextValue $ws.Cells.Item(2, 5) '  +0.03%  '  # E2 (Volume(1h))

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) '1.877.71'  # D3 (Price)
Set-TextValue $ws.Cells.Item(3, 5) '  +0.23%  '  # E3 (Volume(1h))

# Row 4
Set-TextValue $ws.Cells.Item(4, 5) '  +0.20%  '  # E4 (Volume(1h))

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '0.7191'  # D5 (Price)
Set-TextValue $ws.Cells.Item(5, 5) '  +1.25%  '  # E5 (Volume(1h))

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) '242.73'  # D6 (Price)
Set-TextValue $ws.Cells.Item(6, 5) '  +0.40%  '  # E6 (Volume(1h))

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) '1.002'  # D7 (Price)
Set-TextValue $ws.Cells.Item(7, 5) '  +0.15%  '  # E7 (Volume(1h))

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) '0.07995'  # D8 (Price)
Set-TextValue $ws.Cells.Item(8, 5) '  +2.41%  '  # E8 (Volume(1h))

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) '0.3151'  # D9 (Price)
Set-TextValue $ws.Cells.Item(9, 5) '  +1.72%  '  # E9 (Volume(1h))

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) '24.97'  # D10 (Price)
Set-TextValue $ws.Cells.Item(10, 5) '  -0.45%  '  # E10 (Volume(1h))

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) '0.08222'  # D11 (Price)
Set-TextValue $ws.Cells.Item(11, 5) '  -2.14%  '  # E11 (Volume(1h))

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) '1.886.28'  # D12 (Price)
Set-TextValue $ws.Cells.Item(12, 5) '  +0.81%  '  # E12 (Volume(1h))

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) '94.65'  # D13 (Price)
Set-TextValue $ws.Cells.Item(13, 5) '  +3.87%  '  # E13 (Volume(1h))

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) '5.222'  # D14 (Price)
Set-TextValue $ws.Cells.Item(14, 5) '  -0.42%  '  # E14 (Volume(1h))

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) '0.7113'  # D15 (Price)
Set-TextValue $ws.Cells.Item(15, 5) '  -0.01%  '  # E15 (Volume(1h))

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) '6.416'  # D16 (Price)
Set-TextValue $ws.Cells.Item(16, 5) '  +5.64%  '  # E16 (Volume(1h))

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) '0.000008483'  # D17 (Price)
Set-TextValue $ws.Cells.Item(17, 5) '  +3.64%  '  # E17 (Volume(1h))

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) '29.324.44'  # D18 (Price)
Set-TextValue $ws.Cells.Item(18, 5) '  +0.02%  '  # E18 (Volume(1h))

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) '243.20'  # D19 (Price)
Set-TextValue $ws.Cells.Item(19, 5) '  +1.34%  '  # E19 (Volume(1h))

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) '13.27'  # D20 (Price)
Set-TextValue $ws.Cells.Item(20, 5) '  +0.32%  '  # E20 (Volume(1h))

# Row 21
Set-TextValue $ws.Cells.Item(21, 5) '  +0.11%  '  # E21 (Volume(1h))

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) '7.756'  # D22 (Price)
Set-TextValue $ws.Cells.Item(22, 5) '  -0.04%  '  # E22 (Volume(1h))

# Row 23
Set-TextValue $ws.Cells.Item(23, 5) '  +0.14%  '  # E23 (Volume(1h))

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) '0.1595'  # D24 (Price)
Set-TextValue $ws.Cells.Item(24, 5) '  +0.18%  '  # E24 (Volume(1h))

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) '162.53'  # D25 (Price)

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) '9.032'  # D26 (Price)
Set-TextValue $ws.Cells.Item(26, 5) '  +0.38%  '  # E26 (Volume(1h))

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) '18.50'  # D27 (Price)
Set-TextValue $ws.Cells.Item(27, 5) '  +0.22%  '  # E27 (Volume(1h))

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) '1.501'  # D28 (Price)
Set-TextValue $ws.Cells.Item(28, 5) '  -0.23%  '  # E28 (Volume(1h))

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) '4.403'  # D29 (Price)
Set-TextValue $ws.Cells.Item(29, 5) '  +0.25%  '  # E29 (Volume(1h))

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) '4.301'  # D30 (Price)
Set-TextValue $ws.Cells.Item(30, 5) '  +0.15%  '  # E30 (Volume(1h))

# Row 31
Set-TextValue $ws.Cells.Item(31, 5) '  -8.24%  '  # E31 (Volume(1h))

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) '0.05365'  # D32 (Price)
Set-TextValue $ws.Cells.Item(32, 5) '  -0.34%  '  # E32 (Volume(1h))

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) '1.930'  # D33 (Price)
Set-TextValue $ws.Cells.Item(33, 5) '  -0.73%  '  # E33 (Volume(1h))

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) '0.7598'  # D34 (Price)
Set-TextValue $ws.Cells.Item(34, 5) '  +1.21%  '  # E34 (Volume(1h))

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) '1.176'  # D35 (Price)
Set-TextValue $ws.Cells.Item(35, 5) '  -0.12%  '  # E35 (Volume(1h))

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) '2.711'  # D36 (Price)
Set-TextValue $ws.Cells.Item(36, 5) '  +0.63%  '  # E36 (Volume(1h))

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) '1.278.86'  # D37 (Price)
Set-TextValue $ws.Cells.Item(37, 5) '  +3.81%  '  # E37 (Volume(1h))

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) '0.01869'  # D38 (Price)
Set-TextValue $ws.Cells.Item(38, 5) '  -0.03%  '  # E38 (Volume(1h))

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) '2.755'  # D39 (Price)
Set-TextValue $ws.Cells.Item(39, 5) '  +0.93%  '  # E39 (Volume(1h))

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) '6.443'  # D40 (Price)
Set-TextValue $ws.Cells.Item(40, 5) '  -1.15%  '  # E40 (Volume(1h))

# Row 41
Set-TextValue $ws.Cells.Item(41, 2) 'TrustWalletToken'  # B41 (Coin)
Set-TextValue $ws.Cells.Item(41, 3) 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'  # C41 (Link)
Set-TextValue $ws.Cells.Item(41, 4) '0.9134'  # D41 (Price)
Set-TextValue $ws.Cells.Item(41, 5) '  +2.75%  '  # E41 (Volume(1h))

# Row 42
Set-TextValue $ws.Cells.Item(42, 2) 'Quant'  # B42 (Coin)
Set-TextValue $ws.Cells.Item(42, 3) 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'  # C42 (Link)
Set-TextValue $ws.Cells.Item(42, 4) '112.78'  # D42 (Price)
Set-TextValue $ws.Cells.Item(42, 5) '  +3.73%  '  # E42 (Volume(1h))

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '74.13'  # D43 (Price)
Set-TextValue $ws.Cells.Item(43, 5) '  +2.37%  '  # E43 (Volume(1h))

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) '0.00000000131'  # D44 (Price)
Set-TextValue $ws.Cells.Item(44, 5) '  +7.71%  '  # E44 (Volume(1h))

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) '1.002'  # D45 (Price)
Set-TextValue $ws.Cells.Item(45, 5) '  +0.17%  '  # E45 (Volume(1h))

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) '2.026.88'  # D46 (Price)
Set-TextValue $ws.Cells.Item(46, 5) '  +0.26%  '  # E46 (Volume(1h))

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) '0.5230'  # D47 (Price)
Set-TextValue $ws.Cells.Item(47, 5) '  +0.69%  '  # E47 (Volume(1h))

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) '1.793'  # D48 (Price)
Set-TextValue $ws.Cells.Item(48, 5) '  +0.03%  '  # E48 (Volume(1h))

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) '9.488'  # D49 (Price)
Set-TextValue $ws.Cells.Item(49, 5) '  +0.70%  '  # E49 (Volume(1h))

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) '0.4340'  # D50 (Price)
Set-TextValue $ws.Cells.Item(50, 5) '  +0.60%  '  # E50 (Volume(1h))

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) '7.091'  # D51 (Price)
Set-TextValue $ws.Cells.Item(51, 5) '  +0.22%  '  # E51 (Volume(1h))
